# fix for statement body highlighting
# Adds column H ("OK"/"?" status) results for a batch of test rows,
# adds a new "повтор, подсвечивается все тело" test-case row (39),
# bumps row 13's height, and moves the sheet view / selection down
# to where the new row lives.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column H gets "OK" for most of the rows that were missing it.
$okRows = @(3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,22,23,24,25,26,27,28,29,30,31,35,36,37)
foreach ($r in $okRows) {
    $ws.Cells.Item($r, 8).Value = "OK"
}

# Row 21 is a "?" result rather than "OK".
$ws.Cells.Item(21, 8).Value = "?"

# Row 13's header wraps onto an extra line in the edited workbook.
$ws.Rows.Item(13).RowHeight = 18.75

# New test case: "повтор, подсвечивается все тело" (row 39).
$ws.Cells.Item(39, 2).Value = "повтор, подсвечивается все тело"
$ws.Cells.Item(39, 2).WrapText = $true
$ws.Cells.Item(39, 8).Value = "OK"

# Scroll the sheet view down to the newly added row and move the selection.
$ws.Range("H40").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
